$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-fill formatting (date number format for col A, text style for col B) on the new
# rows by copying the existing row 5 formats down before writing values into them.
$ws.Range("A5:B5").Copy()
$ws.Range("A6:B53").PasteSpecial(-4122)

# Carga completa de mayo: full monthly data Apr-2022 .. Apr-2023 (rows 2-53)
$ws.Cells.Item(2, 1).Value = 44652
$ws.Cells.Item(2, 2).Value = "Terapia física"
$ws.Cells.Item(2, 3).Value = 10
$ws.Cells.Item(3, 1).Value = 44652
$ws.Cells.Item(3, 2).Value = "Terapia respiratoria"
$ws.Cells.Item(3, 3).Value = 25
$ws.Cells.Item(4, 1).Value = 44652
$ws.Cells.Item(4, 2).Value = "Terapia ocupacional"
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(5, 1).Value = 44652
$ws.Cells.Item(5, 2).Value = "Terapia de lenguaje"
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(6, 1).Value = 44682
$ws.Cells.Item(6, 2).Value = "Terapia física"
$ws.Cells.Item(6, 3).Value = 75
$ws.Cells.Item(7, 1).Value = 44682
$ws.Cells.Item(7, 2).Value = "Terapia respiratoria"
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(8, 1).Value = 44682
$ws.Cells.Item(8, 2).Value = "Terapia ocupacional"
$ws.Cells.Item(8, 3).Value = 20
$ws.Cells.Item(9, 1).Value = 44682
$ws.Cells.Item(9, 2).Value = "Terapia de lenguaje"
$ws.Cells.Item(9, 3).Value = 20
$ws.Cells.Item(10, 1).Value = 44713
$ws.Cells.Item(10, 2).Value = "Terapia física"
$ws.Cells.Item(10, 3).Value = 171
$ws.Cells.Item(11, 1).Value = 44713
$ws.Cells.Item(11, 2).Value = "Terapia respiratoria"
$ws.Cells.Item(11, 3).Value = 13
$ws.Cells.Item(12, 1).Value = 44713
$ws.Cells.Item(12, 2).Value = "Terapia ocupacional"
$ws.Cells.Item(12, 3).Value = 20
$ws.Cells.Item(13, 1).Value = 44713
$ws.Cells.Item(13, 2).Value = "Terapia de lenguaje"
$ws.Cells.Item(13, 3).Value = 20
$ws.Cells.Item(14, 1).Value = 44743
$ws.Cells.Item(14, 2).Value = "Terapia física"
$ws.Cells.Item(14, 3).Value = 319
$ws.Cells.Item(15, 1).Value = 44743
$ws.Cells.Item(15, 2).Value = "Terapia respiratoria"
$ws.Cells.Item(15, 3).Value = 56
$ws.Cells.Item(16, 1).Value = 44743
$ws.Cells.Item(16, 2).Value = "Terapia ocupacional"
$ws.Cells.Item(16, 3).Value = 31
$ws.Cells.Item(17, 1).Value = 44743
$ws.Cells.Item(17, 2).Value = "Terapia de lenguaje"
$ws.Cells.Item(17, 3).Value = 14
$ws.Cells.Item(18, 1).Value = 44774
$ws.Cells.Item(18, 2).Value = "Terapia física"
$ws.Cells.Item(18, 3).Value = 1268
$ws.Cells.Item(19, 1).Value = 44774
$ws.Cells.Item(19, 2).Value = "Terapia respiratoria"
$ws.Cells.Item(19, 3).Value = 31
$ws.Cells.Item(20, 1).Value = 44774
$ws.Cells.Item(20, 2).Value = "Terapia ocupacional"
$ws.Cells.Item(20, 3).Value = 117
$ws.Cells.Item(21, 1).Value = 44774
$ws.Cells.Item(21, 2).Value = "Terapia de lenguaje"
$ws.Cells.Item(21, 3).Value = 147
$ws.Cells.Item(22, 1).Value = 44805
$ws.Cells.Item(22, 2).Value = "Terapia física"
$ws.Cells.Item(22, 3).Value = 1480
$ws.Cells.Item(23, 1).Value = 44805
$ws.Cells.Item(23, 2).Value = "Terapia respiratoria"
$ws.Cells.Item(23, 3).Value = 62
$ws.Cells.Item(24, 1).Value = 44805
$ws.Cells.Item(24, 2).Value = "Terapia ocupacional"
$ws.Cells.Item(24, 3).Value = 283
$ws.Cells.Item(25, 1).Value = 44805
$ws.Cells.Item(25, 2).Value = "Terapia de lenguaje"
$ws.Cells.Item(25, 3).Value = 79
$ws.Cells.Item(26, 1).Value = 44835
$ws.Cells.Item(26, 2).Value = "Terapia física"
$ws.Cells.Item(26, 3).Value = 1454
$ws.Cells.Item(27, 1).Value = 44835
$ws.Cells.Item(27, 2).Value = "Terapia respiratoria"
$ws.Cells.Item(27, 3).Value = 47
$ws.Cells.Item(28, 1).Value = 44835
$ws.Cells.Item(28, 2).Value = "Terapia ocupacional"
$ws.Cells.Item(28, 3).Value = 341
$ws.Cells.Item(29, 1).Value = 44835
$ws.Cells.Item(29, 2).Value = "Terapia de lenguaje"
$ws.Cells.Item(29, 3).Value = 94
$ws.Cells.Item(30, 1).Value = 44866
$ws.Cells.Item(30, 2).Value = "Terapia física"
$ws.Cells.Item(30, 3).Value = 1426
$ws.Cells.Item(31, 1).Value = 44866
$ws.Cells.Item(31, 2).Value = "Terapia respiratoria"
$ws.Cells.Item(31, 3).Value = 66
$ws.Cells.Item(32, 1).Value = 44866
$ws.Cells.Item(32, 2).Value = "Terapia ocupacional"
$ws.Cells.Item(32, 3).Value = 311
$ws.Cells.Item(33, 1).Value = 44866
$ws.Cells.Item(33, 2).Value = "Terapia de lenguaje"
$ws.Cells.Item(33, 3).Value = 123
$ws.Cells.Item(34, 1).Value = 44896
$ws.Cells.Item(34, 2).Value = "Terapia física"
$ws.Cells.Item(34, 3).Value = 1418
$ws.Cells.Item(35, 1).Value = 44896
$ws.Cells.Item(35, 2).Value = "Terapia respiratoria"
$ws.Cells.Item(35, 3).Value = 90
$ws.Cells.Item(36, 1).Value = 44896
$ws.Cells.Item(36, 2).Value = "Terapia ocupacional"
$ws.Cells.Item(36, 3).Value = 253
$ws.Cells.Item(37, 1).Value = 44896
$ws.Cells.Item(37, 2).Value = "Terapia de lenguaje"
$ws.Cells.Item(37, 3).Value = 133
$ws.Cells.Item(38, 1).Value = 44927
$ws.Cells.Item(38, 2).Value = "Terapia física"
$ws.Cells.Item(38, 3).Value = 1098
$ws.Cells.Item(39, 1).Value = 44927
$ws.Cells.Item(39, 2).Value = "Terapia respiratoria"
$ws.Cells.Item(39, 3).Value = 76
$ws.Cells.Item(40, 1).Value = 44927
$ws.Cells.Item(40, 2).Value = "Terapia ocupacional"
$ws.Cells.Item(40, 3).Value = 164
$ws.Cells.Item(41, 1).Value = 44927
$ws.Cells.Item(41, 2).Value = "Terapia de lenguaje"
$ws.Cells.Item(41, 3).Value = 68
$ws.Cells.Item(42, 1).Value = 44958
$ws.Cells.Item(42, 2).Value = "Terapia física"
$ws.Cells.Item(42, 3).Value = 1238
$ws.Cells.Item(43, 1).Value = 44958
$ws.Cells.Item(43, 2).Value = "Terapia respiratoria"
$ws.Cells.Item(43, 3).Value = 84
$ws.Cells.Item(44, 1).Value = 44958
$ws.Cells.Item(44, 2).Value = "Terapia ocupacional"
$ws.Cells.Item(44, 3).Value = 248
$ws.Cells.Item(45, 1).Value = 44958
$ws.Cells.Item(45, 2).Value = "Terapia de lenguaje"
$ws.Cells.Item(45, 3).Value = 73
$ws.Cells.Item(46, 1).Value = 44986
$ws.Cells.Item(46, 2).Value = "Terapia física"
$ws.Cells.Item(46, 3).Value = 1356
$ws.Cells.Item(47, 1).Value = 44986
$ws.Cells.Item(47, 2).Value = "Terapia respiratoria"
$ws.Cells.Item(47, 3).Value = 92
$ws.Cells.Item(48, 1).Value = 44986
$ws.Cells.Item(48, 2).Value = "Terapia ocupacional"
$ws.Cells.Item(48, 3).Value = 296
$ws.Cells.Item(49, 1).Value = 44986
$ws.Cells.Item(49, 2).Value = "Terapia de lenguaje"
$ws.Cells.Item(49, 3).Value = 152
$ws.Cells.Item(50, 1).Value = 45017
$ws.Cells.Item(50, 2).Value = "Terapia física"
$ws.Cells.Item(50, 3).Value = 1127
$ws.Cells.Item(51, 1).Value = 45017
$ws.Cells.Item(51, 2).Value = "Terapia respiratoria"
$ws.Cells.Item(51, 3).Value = 40
$ws.Cells.Item(52, 1).Value = 45017
$ws.Cells.Item(52, 2).Value = "Terapia ocupacional"
$ws.Cells.Item(52, 3).Value = 133
$ws.Cells.Item(53, 1).Value = 45017
$ws.Cells.Item(53, 2).Value = "Terapia de lenguaje"
$ws.Cells.Item(53, 3).Value = 94

# Turn on the AutoFilter over the whole table and register the hidden filter-database name
$ws.Range("A1:C53").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", $ws.Range("A1:C53"))
$filterName.Visible = $false

# Leave the selection where the author left it when saving
$ws.Range("B1").Select()
